$wb = $excel.ActiveWorkbook

# This script applies updated market-price figures (columns H-N) that were
# refreshed by the scheduled data-pull runner, sheet by sheet.

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 931.44446
$ws.Range("I4").Value = 797.875
$ws.Range("K4").Value = 797.875
$ws.Range("M4").Value = -683.875
$ws.Range("H17").Value = 2010.3636
$ws.Range("J17").Value = 2010.3636
$ws.Range("L17").Value = 6031.0908
$ws.Range("N17").Value = -6367.0908
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1087.4706
$ws.Range("J2").Value = 1357.8334
$ws.Range("L2").Value = 1357.8334
$ws.Range("N2").Value = -1583.8334
$ws.Range("H32").Value = 3642.8936
$ws.Range("I32").Value = 3642.8936
$ws.Range("K32").Value = 3642.8936
$ws.Range("M32").Value = -3355.8936
$ws.Range("H74").Value = 207532.6
$ws.Range("I74").Value = 278888.4
$ws.Range("J74").Value = 3658.8572
$ws.Range("K74").Value = 278888.4
$ws.Range("L74").Value = 3658.8572
$ws.Range("M74").Value = -278014.4
$ws.Range("N74").Value = -5406.8572
$ws.Range("H77").Value = 207532.6
$ws.Range("I77").Value = 278888.4
$ws.Range("J77").Value = 3658.8572
$ws.Range("K77").Value = 1394442
$ws.Range("L77").Value = 18294.286
$ws.Range("M77").Value = -1390074
$ws.Range("N77").Value = -27030.286
$ws.Range("H88").Value = 4171.222
$ws.Range("I88").Value = 2058.25
$ws.Range("K88").Value = 2058.25
$ws.Range("M88").Value = -1652.25
$ws.Range("H91").Value = 4171.222
$ws.Range("I91").Value = 2058.25
$ws.Range("K91").Value = 2058.25
$ws.Range("M91").Value = -654.25
$ws.Range("H116").Value = 1087.4706
$ws.Range("J116").Value = 1357.8334
$ws.Range("L116").Value = 1357.8334
$ws.Range("N116").Value = -5945.8334

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1087.4706
$ws.Range("J3").Value = 1357.8334
$ws.Range("L3").Value = 1357.8334
$ws.Range("N3").Value = -1585.8334
$ws.Range("H86").Value = 3379.2856
$ws.Range("I86").Value = 3379.2856
$ws.Range("K86").Value = 3379.2856
$ws.Range("M86").Value = -2256.2856
$ws.Range("H89").Value = 3379.2856
$ws.Range("I89").Value = 3379.2856
$ws.Range("K89").Value = 16896.428
$ws.Range("M89").Value = -11280.428
$ws.Range("H105").Value = 13002468
$ws.Range("I105").Value = 770866.25
$ws.Range("J105").Value = 35718300
$ws.Range("K105").Value = 770866.25
$ws.Range("L105").Value = 35718300
$ws.Range("M105").Value = -769119.25
$ws.Range("N105").Value = -35721794
$ws.Range("H134").Value = 5750
$ws.Range("I134").Value = 1250
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 3750
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -1215
$ws.Range("N134").Value = -29070

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4033.575
$ws.Range("I31").Value = 2689.1614
$ws.Range("J31").Value = 8664.333000000001
$ws.Range("K31").Value = 2689.1614
$ws.Range("L31").Value = 8664.333000000001
$ws.Range("M31").Value = -2394.1614
$ws.Range("N31").Value = -9254.333000000001
$ws.Range("H34").Value = 4033.575
$ws.Range("I34").Value = 2689.1614
$ws.Range("J34").Value = 8664.333000000001
$ws.Range("K34").Value = 2689.1614
$ws.Range("L34").Value = 8664.333000000001
$ws.Range("M34").Value = -2487.1614
$ws.Range("N34").Value = -9068.333000000001
$ws.Range("H57").Value = 31185
$ws.Range("J57").Value = 21750
$ws.Range("L57").Value = 21750
$ws.Range("N57").Value = -22870
$ws.Range("H122").Value = 3858.3333
$ws.Range("J122").Value = 5592
$ws.Range("L122").Value = 16776
$ws.Range("N122").Value = -21676
$ws.Range("H132").Value = 3731.5
$ws.Range("I132").Value = 3733.1333
$ws.Range("J132").Value = 3728.7778
$ws.Range("K132").Value = 11199.3999
$ws.Range("L132").Value = 11186.3334
$ws.Range("M132").Value = -8669.3999
$ws.Range("N132").Value = -16246.3334
$ws.Range("H134").Value = 3636.1072
$ws.Range("I134").Value = 3367.2083
$ws.Range("J134").Value = 5249.5
$ws.Range("K134").Value = 10101.6249
$ws.Range("L134").Value = 15748.5
$ws.Range("M134").Value = -7566.624899999999
$ws.Range("N134").Value = -20818.5
$ws.Range("H141").Value = 655093.3
$ws.Range("I141").Value = 60296
$ws.Range("J141").Value = 721181.9
$ws.Range("K141").Value = 60296
$ws.Range("L141").Value = 721181.9
$ws.Range("M141").Value = -55116
$ws.Range("N141").Value = -731541.9

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 2518
$ws.Range("I24").Value = 1748.75
$ws.Range("K24").Value = 5246.25
$ws.Range("M24").Value = -5016.25
$ws.Range("H116").Value = 2182998.8
$ws.Range("J116").Value = 4997.5
$ws.Range("L116").Value = 14992.5
$ws.Range("N116").Value = -21876.5
$ws.Range("H122").Value = 1656.6666
$ws.Range("J122").Value = 2005
$ws.Range("L122").Value = 18045
$ws.Range("N122").Value = -22945
$ws.Range("H130").Value = 10520
$ws.Range("J130").Value = 6500
$ws.Range("L130").Value = 19500
$ws.Range("N130").Value = -29540
$ws.Range("H139").Value = 1910.0555
$ws.Range("I139").Value = 1370.0714
$ws.Range("K139").Value = 4110.2142
$ws.Range("M139").Value = 1029.7858
$ws.Range("H140").Value = 5601.8
$ws.Range("I140").Value = 5371.522
$ws.Range("J140").Value = 8250
$ws.Range("K140").Value = 16114.566
$ws.Range("L140").Value = 24750
$ws.Range("M140").Value = -10934.566
$ws.Range("N140").Value = -35110

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 657.53845
$ws.Range("I2").Value = 858.44446
$ws.Range("K2").Value = 858.44446
$ws.Range("M2").Value = -745.44446
$ws.Range("H126").Value = 11331.143
$ws.Range("J126").Value = 17998
$ws.Range("L126").Value = 53994
$ws.Range("N126").Value = -58934
$ws.Range("H132").Value = 7747.9287
$ws.Range("I132").Value = 2121.625
$ws.Range("K132").Value = 6364.875
$ws.Range("M132").Value = -3834.875
$ws.Range("H141").Value = 73165.664
$ws.Range("J141").Value = 73165.664
$ws.Range("L141").Value = 73165.664
$ws.Range("N141").Value = -83525.664

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4616.0713
$ws.Range("I40").Value = 4621.8887
$ws.Range("K40").Value = 4621.8887
$ws.Range("M40").Value = -4485.8887
$ws.Range("H46").Value = 1360.625
$ws.Range("I46").Value = 647.6667
$ws.Range("J46").Value = 1788.4
$ws.Range("K46").Value = 647.6667
$ws.Range("L46").Value = 1788.4
$ws.Range("M46").Value = -459.6667
$ws.Range("N46").Value = -2164.4
$ws.Range("H61").Value = 3390.3572
$ws.Range("I61").Value = 3553.25
$ws.Range("K61").Value = 3553.25
$ws.Range("M61").Value = -3351.25
$ws.Range("H113").Value = 3390.3572
$ws.Range("I113").Value = 3553.25
$ws.Range("K113").Value = 3553.25
$ws.Range("M113").Value = -1383.25
$ws.Range("H122").Value = 2693.4285
$ws.Range("I122").Value = 2695.9
$ws.Range("K122").Value = 8087.700000000001
$ws.Range("M122").Value = -5637.700000000001
$ws.Range("H136").Value = 3907.4827
$ws.Range("I136").Value = 3692.1428
$ws.Range("K136").Value = 11076.4284
$ws.Range("M136").Value = -8526.428400000001

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 17588.8
$ws.Range("I14").Value = 15678.615
$ws.Range("K14").Value = 15678.615
$ws.Range("M14").Value = -15510.615
$ws.Range("H70").Value = 28994
$ws.Range("J70").Value = 28994
$ws.Range("L70").Value = 28994
$ws.Range("N70").Value = -29624
$ws.Range("H73").Value = 28994
$ws.Range("J73").Value = 28994
$ws.Range("L73").Value = 28994
$ws.Range("N73").Value = -31178
$ws.Range("H136").Value = 58827812
$ws.Range("I136").Value = 66667856
$ws.Range("K136").Value = 200003568
$ws.Range("M136").Value = -200001018
$ws.Range("H140").Value = 148973.5
$ws.Range("J140").Value = 148973.5
$ws.Range("L140").Value = 148973.5
$ws.Range("N140").Value = -159333.5
$ws.Range("H141").Value = 70036.25
$ws.Range("J141").Value = 70036.25
$ws.Range("L141").Value = 70036.25
$ws.Range("N141").Value = -80396.25
